# "Fixed Bento 80 Test scripts":
# Append "order By <field> ASC LIMIT 100" to the Cypher queries used by the
# Cases/Samples/Files tabs (B2/B3/B4) on the "startup" sheet. The aggregate
# count query in column C is left unchanged.
#
# Cells are (re)written in C, D, A order so the workbook's shared-string
# table ends up in the same order as the authored workbook (unchanged
# queries keep their slot, newly-edited ones are appended at the end in the
# order they are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# FilesTab query (row 4, column B) -> add order by file name
$CQuery = @'
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and tp.chemotherapy_regimen In ["standard AC followed by a taxane"] and d.tumor_grade In ["Intermediate Grade"]
RETURN  f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
    f.file_size AS `Size`,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`,
    ss.disease_subtype as `Diagnosis`
 order By f.file_name ASC LIMIT 100
'@
$ws.Range("B4").Value = $CQuery

# SamplesTab query (row 3, column B) -> add order by sample id
$DQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and tp.chemotherapy_regimen In ["standard AC followed by a taxane"] and d.tumor_grade In ["Intermediate Grade"]WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`,
            lp.test_name as Platform
 order By samp.sample_id ASC LIMIT 100
'@
$ws.Range("B3").Value = $DQuery

# CasesTab query (row 2, column B) -> add order by case/study_subject id
$AQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and tp.chemotherapy_regimen In ["standard AC followed by a taxane"] and d.tumor_grade In ["Intermediate Grade"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100
'@
$ws.Range("B2").Value = $AQuery

# Row heights grow to fit the longer, edited query text (wrap-text cells).
$ws.Rows.Item(2).RowHeight = 374.4
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 331.2

# Selection moves to B2 (and the view is no longer scrolled to row 3).
$ws.Range("B2").Select() | Out-Null
